$d = $word.ActiveDocument
$tab = [char]9

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (it sat after "GPA: 3.63" and before
#    "/4.0"). Word will recompact bookmark ids on save so the remaining
#    "_30j0zll" bookmark naturally becomes id 1.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# ---------------------------------------------------------------------------
# 2. In the "Computer ... Git ... PHP" line, drop the "Git" entry (it is
#    being relocated to the Skills line below), collapsing the surrounding
#    double-tabs down to a single tab before "PHP".
# ---------------------------------------------------------------------------
$computerPara = $d.Paragraphs.Item(43)
$gitRng = $computerPara.Range.Duplicate
$found = $gitRng.Find.Execute("Git", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $delRng = $d.Range($gitRng.Start - 2, $gitRng.End + 2)
    $delRng.Text = $tab
}

# ---------------------------------------------------------------------------
# 3. In the "Skills ... Code::Blocks  Eclipse  Django" line, remove the two
#    editors "Code::Blocks" and "Eclipse" (and their separating tabs), then
#    append "Git" (moved from the Computer line above) after "Django".
# ---------------------------------------------------------------------------
$skillsPara = $d.Paragraphs.Item(44)
$codeRng = $skillsPara.Range.Duplicate
$foundCode = $codeRng.Find.Execute("Code::Blocks", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundCode) {
    $djangoRng = $skillsPara.Range.Duplicate
    $djangoRng.Find.Execute("Django", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $delRng2 = $d.Range($codeRng.Start, $djangoRng.Start)
    $delRng2.Text = ""
}

$skillsPara2 = $d.Paragraphs.Item(44)
$djangoRng2 = $skillsPara2.Range.Duplicate
$djangoRng2.Find.Execute("Django", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint = $d.Range($djangoRng2.End, $djangoRng2.End)
$insPoint.Text = "$tab$tab" + "Git"

# ---------------------------------------------------------------------------
# 4. Re-create "_GoBack" at the new edit location, right after the "Git"
#    that now ends the Skills line (this mirrors where Word leaves the
#    mark after the last edit made to the document).
# ---------------------------------------------------------------------------
$skillsPara3 = $d.Paragraphs.Item(44)
$gitRng2 = $skillsPara3.Range.Duplicate
$gitRng2.Find.Execute("Git", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRng = $d.Range($gitRng2.End, $gitRng2.End)
$goBackRng.Bookmarks.Add("_GoBack")
